$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A42").Value = "GRT-USD"
